$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp note in A1
$ws.Range("A1").Value = "Datos actualizados a 4 de Abril de 2020 a las 16:52"

# Country data table updates: a handful of countries (Moldavia, Bielorrusia,
# Gibraltar, Zambia, Haiti) moved to new rows (alphabetising/resorting a block
# of the country list) and the case-count columns were refreshed with newer
# figures, which cascades a handful of neighbouring rows down by one.
$rows = @{}
$rows[4] = @("Estados Unidos", 279500, 2339, 12729, 259314, 5804, 65, 7457)
$rows[7] = @("Alemania", 92081, 922, 26400, 64382, 3936, 24, 1299)
$rows[44] = @("Panama", 1673, 0, 13, 1619, 50, 0, 41)
$rows[54] = @("Singapur", 1189, 75, 297, 886, 24, 1, 6)
$rows[68] = @("Moldavia", 752, 161, 26, 716, 65, 2, 10)
$rows[69] = @("Crucero", 712, 0, 619, 82, 10, 0, 11)
$rows[70] = @("Barein", 688, 16, 399, 285, 3, 0, 4)
$rows[71] = @("Hungria", 678, 55, 58, 588, 17, 6, 32)
$rows[72] = @("Bosnia y Herzegovina", 617, 38, 28, 570, 4, 2, 19)
$rows[78] = @("Bulgaria", 503, 18, 34, 452, 26, 3, 17)
$rows[84] = @("Bielorrusia", 440, 89, 53, 382, 11, 1, 5)
$rows[85] = @("Costa Rica", 416, 0, 11, 403, 13, 0, 2)
$rows[86] = @("Republica de Chipre", 396, 0, 28, 357, 11, 0, 11)
$rows[123] = @("Trinidad yTobago", 101, 3, 1, 94, 0, 0, 6)
$rows[124] = @("Gibraltar", 98, 3, 52, 46, 0, 0, 0)
$rows[125] = @("Paraguay", 96, 4, 12, 81, 2, 0, 3)
$rows[144] = @("Zambia", 39, 0, 2, 36, 0, 0, 1)
$rows[145] = @("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
$rows[157] = @("Haiti", 20, 2, 1, 19, 0, 0, 0)
$rows[158] = @("Birmania", 20, 0, 0, 19, 0, 0, 1)
$rows[159] = @("Tanzania", 20, 0, 3, 16, 0, 0, 1)
$rows[160] = @("Maldivas", 19, 0, 13, 6, 0, 0, 0)
$rows[161] = @("Nueva Caledonia", 18, 0, 1, 17, 0, 0, 0)

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($c = 1; $c -le 8; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}
